$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H12" = 440
    "I12" = 127.5
    "J12" = 690
    "K12" = 127.5
    "L12" = 690
    "M12" = 42.5
    "N12" = -1030
    "H17" = 2345.7307
    "J17" = 2666.1428
    "L17" = 7998.428400000001
    "N17" = -8334.428400000001
    "H43" = 1800.2
    "I43" = 2000.5
    "J43" = 1666.6666
    "K43" = 2000.5
    "L43" = 1666.6666
    "M43" = -1931.5
    "N43" = -1804.6666
    "H80" = 1281
    "I80" = 450
    "K80" = 1350
    "M80" = -352
    "H83" = 1281
    "I83" = 450
    "K83" = 4050
    "M83" = 942
    "H92" = 165.77777
    "I92" = 172.83333
    "J92" = 151.66667
    "K92" = 172.83333
    "L92" = 151.66667
    "M92" = 1075.16667
    "N92" = -2647.66667
    "H98" = 920.5
    "I98" = 737.1111
    "J98" = 1470.6666
    "K98" = 737.1111
    "L98" = 1470.6666
    "M98" = 760.8889
    "N98" = -4466.6666
    "H122" = 920.5
    "I122" = 737.1111
    "J122" = 1470.6666
    "K122" = 2211.3333
    "L122" = 4411.9998
    "M122" = 238.6667000000002
    "N122" = -9311.9998
    "H124" = 0
    "J124" = 0
    "N124" = 0
    "H138" = 3009.7407
    "J138" = 3627.4707
    "L138" = 10882.4121
    "N138" = -21162.4121
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("L124").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H74" = 5297.5454
    "I74" = 5196.96
    "J74" = 5611.875
    "K74" = 5196.96
    "L74" = 5611.875
    "M74" = -4322.96
    "N74" = -7359.875
    "H77" = 5297.5454
    "I77" = 5196.96
    "J77" = 5611.875
    "K77" = 25984.8
    "L77" = 28059.375
    "M77" = -21616.8
    "N77" = -36795.375
    "H98" = 10000
    "J98" = 10000
    "L98" = 10000
    "N98" = -15990
    "H132" = 2577.4
    "J132" = 0
    "L132" = 0
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N132").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H86" = 2869.35
    "I86" = 1158.6666
    "K86" = 1158.6666
    "M86" = -35.66660000000002
    "H88" = 1136181.1
    "J88" = 1136181.1
    "L88" = 1136181.1
    "N88" = -1136993.1
    "H89" = 2869.35
    "I89" = 1158.6666
    "K89" = 5793.333000000001
    "M89" = -177.3330000000005
    "H91" = 1136181.1
    "J91" = 1136181.1
    "L91" = 1136181.1
    "N91" = -1138989.1
    "H95" = 7384
    "J95" = 7384
    "L95" = 7384
    "N95" = -12876
    "H100" = 10528.5
    "J100" = 10528.5
    "L100" = 10528.5
    "N100" = -12692.5
    "H105" = 735.7143
    "I105" = 735.7143
    "K105" = 735.7143
    "M105" = 1011.2857
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H22" = 1673.25
    "I22" = 663.55554
    "K22" = 663.55554
    "M22" = -313.55554
    "H25" = 1737.375
    "I25" = 1737.375
    "K25" = 1737.375
    "M25" = -1563.375
    "H28" = 16970.5
    "J28" = 16970.5
    "L28" = 16970.5
    "N28" = -17460.5
    "H58" = 4294.143
    "I58" = 2243.25
    "J58" = 7028.6665
    "K58" = 2243.25
    "L58" = 7028.6665
    "M58" = -2040.25
    "N58" = -7434.6665
    "H87" = 0
    "I87" = 0
    "K87" = 0
    "H88" = 8085.25
    "J88" = 8085.25
    "L88" = 8085.25
    "N88" = -8897.25
    "H90" = 0
    "I90" = 0
    "K90" = 0
    "H91" = 8085.25
    "J91" = 8085.25
    "L91" = 8085.25
    "N91" = -10893.25
    "H132" = 4336.1577
    "I132" = 3997.125
    "K132" = 11991.375
    "M132" = -9461.375
    "H134" = 2279
    "I134" = 2351.375
    "K134" = 7054.125
    "M134" = -4519.125
    "H136" = 4294.143
    "I136" = 2243.25
    "J136" = 7028.6665
    "K136" = 6729.75
    "L136" = 21085.9995
    "M136" = -4179.75
    "N136" = -26185.9995
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("M87").ClearContents()
$ws.Range("M90").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H139" = 5027.4287
    "I139" = 4602.8335
    "K139" = 13808.5005
    "M139" = -8668.500499999998
    "H140" = 4814.25
    "I140" = 4669
    "K140" = 14007
    "M140" = -8827
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H52" = 41000
    "J52" = 41000
    "L52" = 41000
    "N52" = -41518
    "H70" = 3079.0833
    "J70" = 0
    "L70" = 0
    "H73" = 3079.0833
    "J73" = 0
    "L73" = 0
    "H80" = 3212.5
    "I80" = 3533.3333
    "J80" = 2250
    "K80" = 3533.3333
    "L80" = 2250
    "M80" = -2535.3333
    "N80" = -4246
    "H83" = 3212.5
    "I83" = 3533.3333
    "J83" = 2250
    "K83" = 17666.6665
    "L83" = 11250
    "M83" = -12674.6665
    "N83" = -21234
    "H122" = 3715.889
    "I122" = 3135.4285
    "K122" = 9406.2855
    "M122" = -6956.2855
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H40" = 5232.8
    "I40" = 5721.3335
    "J40" = 4500
    "K40" = 5721.3335
    "L40" = 4500
    "M40" = -5585.3335
    "N40" = -4772
    "H55" = 1046.6
    "I55" = 2035
    "J55" = 387.66666
    "K55" = 2035
    "L55" = 387.66666
    "M55" = -1862
    "N55" = -733.66666
    "H61" = 3873.0667
    "I61" = 1788.6666
    "K61" = 1788.6666
    "M61" = -1586.6666
    "H113" = 3873.0667
    "I113" = 1788.6666
    "K113" = 1788.6666
    "M113" = 381.3334
    "H132" = 6947.2856
    "I132" = 2656.5
    "K132" = 7969.5
    "M132" = -5439.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H2" = 1230.5555
    "I2" = 1374.375
    "J2" = 80
    "K2" = 1374.375
    "L2" = 80
    "M2" = -1262.375
    "N2" = -304
    "H4" = 1856268.1
    "I4" = 1856268.1
    "K4" = 1856268.1
    "M4" = -1856155.1
    "H122" = 2500
    "I122" = 2500
    "K122" = 7500
    "M122" = -5050
    "H132" = 1626.9412
    "I132" = 1212.1538
    "K132" = 3636.4614
    "M132" = -1106.4614
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
